$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric card-number string that must stay plain text
# (same as the original cell). Assigning a digit-only string directly would
# make Excel auto-convert it to a number, so enter it with a leading
# apostrophe (forces text / quotePrefix) and then copy the number
# format/style from the still-text neighbour C3 back onto B3 so the cell's
# style index is unaffected by the quote-prefix style bump.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 23.01.2024"

$ws.Range("B6").Value = "27.01."
$ws.Range("C6").Value = "28.01."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 65769082"
$ws.Range("E6").Value = "39,36-"

$ws.Range("B7").Value = "30.01."
$ws.Range("C7").Value = "31.01."
$ws.Range("D7").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E7").Value = "76,80-"

$ws.Range("B8").Value = "01.02."
$ws.Range("C8").Value = "02.02."
$ws.Range("D8").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E8").Value = "61,71-"

$ws.Range("B9").Value = "03.02."
$ws.Range("C9").Value = "04.02."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "25,18-"

$ws.Range("B10").Value = "04.02."
$ws.Range("C10").Value = "05.02."
$ws.Range("D10").Value = "MCDONALDS Badibling"
$ws.Range("E10").Value = "35,48-"

$ws.Range("D12").Value = "KONTOSTAND AM 09.02.2024"
$ws.Range("E12").Value = "238,53-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 16.02.2024"
